# Update factsheets with text edits from COMM
#
# The source data previously stored the "No. of 990 Filers w/ Gov Grants"
# counts as real numbers. This edit converts every one of those numeric
# cells into a plain text cell holding the same digits (to match the
# text-based formatting used by the rest of the sheet), fixes up the
# previously-blank "Pleasants County" row on the County sheet, and adds a
# new "Total" row underneath it.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    # Drop back to the default (unstyled) cell style so we don't leave a
    # stray "@" number-format behind on the cell.
    $Cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# "Overall" sheet: A2 (total filer count) becomes text.
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextCell $wsOverall.Cells.Item(2, 1) "769"

# ---------------------------------------------------------------------
# "County" sheet: column B (filer count) becomes text for every county
# row (2-54), the previously all-zero "Pleasants County" row (55) gets
# its real percentages/dollar figures, and a new "Total" row (56) is
# appended.
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2 = "8"; 3 = "28"; 4 = "6"; 5 = "4"; 6 = "6"; 7 = "53"; 8 = "3";
    9 = "4"; 10 = "5"; 11 = "1"; 12 = "5"; 13 = "25"; 14 = "8"; 15 = "13";
    16 = "13"; 17 = "24"; 18 = "6"; 19 = "21"; 20 = "111"; 21 = "10";
    22 = "6"; 23 = "13"; 24 = "22"; 25 = "9"; 26 = "6"; 27 = "7"; 28 = "21";
    29 = "14"; 30 = "10"; 31 = "47"; 32 = "4"; 33 = "8"; 34 = "12";
    35 = "44"; 36 = "8"; 37 = "5"; 38 = "15"; 39 = "9"; 40 = "27";
    41 = "20"; 42 = "5"; 43 = "6"; 44 = "4"; 45 = "5"; 46 = "7"; 47 = "3";
    48 = "12"; 49 = "6"; 50 = "2"; 51 = "9"; 52 = "2"; 53 = "41"; 54 = "6"
}

foreach ($r in $countyCounts.Keys) {
    Set-TextCell $wsCounty.Cells.Item($r, 2) $countyCounts[$r]
}

# Row 55 ("Pleasants County") previously held placeholder zeros; fill in
# the real figures.
Set-TextCell $wsCounty.Cells.Item(55, 2) "0.00%"
Set-TextCell $wsCounty.Cells.Item(55, 3) "`$0"
Set-TextCell $wsCounty.Cells.Item(55, 4) "0.00%"
Set-TextCell $wsCounty.Cells.Item(55, 5) "0.00%"
Set-TextCell $wsCounty.Cells.Item(55, 6) "0.00%"

# New row 56: state-wide "Total" row.
Set-TextCell $wsCounty.Cells.Item(56, 1) "Total"
Set-TextCell $wsCounty.Cells.Item(56, 2) "769"
Set-TextCell $wsCounty.Cells.Item(56, 3) "`$1,067,242,685"
Set-TextCell $wsCounty.Cells.Item(56, 4) "9.05%"
Set-TextCell $wsCounty.Cells.Item(56, 5) "-27.48%"
Set-TextCell $wsCounty.Cells.Item(56, 6) "75.81%"

# ---------------------------------------------------------------------
# "Congressional District" sheet: column B becomes text (rows 2-4).
# ---------------------------------------------------------------------
$wsDistrict = $wb.Worksheets.Item("Congressional District")
$districtCounts = @{ 2 = "368"; 3 = "401"; 4 = "769" }
foreach ($r in $districtCounts.Keys) {
    Set-TextCell $wsDistrict.Cells.Item($r, 2) $districtCounts[$r]
}

# ---------------------------------------------------------------------
# "Size" sheet: column B becomes text (rows 2-8).
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @{ 2 = "278"; 3 = "183"; 4 = "111"; 5 = "32"; 6 = "94"; 7 = "71"; 8 = "769" }
foreach ($r in $sizeCounts.Keys) {
    Set-TextCell $wsSize.Cells.Item($r, 2) $sizeCounts[$r]
}

# ---------------------------------------------------------------------
# "Subsector" sheet: column B becomes text (rows 2-12).
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
$subsectorCounts = @{
    2 = "38"; 3 = "32"; 4 = "31"; 5 = "63"; 6 = "13"; 7 = "316";
    8 = "55"; 9 = "7"; 10 = "206"; 11 = "8"; 12 = "769"
}
foreach ($r in $subsectorCounts.Keys) {
    Set-TextCell $wsSubsector.Cells.Item($r, 2) $subsectorCounts[$r]
}
